# Resultados_de_Licitacoes_e_Conjecturas_v2.xlsx — "Add files via upload"
#
# The underlying edit swaps the "Mediana_pre" / "Media_saneada_pre" columns
# (F and G) -- both the header text and every data row -- and nudges a few
# column widths plus the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns F (6) and G (7): header (row 1) + all data rows (2-17) ---
for ($r = 1; $r -le 17; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 6).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $fVal
}

# --- Column width tweaks ---
$ws.Range("B1").ColumnWidth = 14.0
$ws.Range("C1").ColumnWidth = 20.8
$ws.Range("E1").ColumnWidth = 11.3
$ws.Range("F1").ColumnWidth = 11.0
$ws.Range("G1").ColumnWidth = 11.0

# --- Move the saved selection/active cell ---
$ws.Range("G19").Select()
